# Auto-generated: apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.919.75"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.452.78"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'578.81"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "'165.63"
$ws.Range("E6").Value = "  -5.26%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").Value = "2.451.27"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "'4.87"
$ws.Range("E13").Value = "  -4.75%  "
$ws.Range("D14").Value = "'25.33"
$ws.Range("E14").Value = "  -4.78%  "
$ws.Range("D15").Value = "2.888.22"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "66.663.05"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("E17").Value = "  -6.06%  "
$ws.Range("D18").Value = "2.447.79"
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").Value = "'11.33"
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("D20").Value = "'7.73"
$ws.Range("E20").Value = "  -4.83%  "
$ws.Range("D21").Value = "'354.37"
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'69.50"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").Value = "'4.21"
$ws.Range("E25").Value = "  -9.64%  "
$ws.Range("E26").Value = "  -8.59%  "
$ws.Range("D27").Value = "'8.92"
$ws.Range("E27").Value = "  -11.17%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "2.570.53"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").Value = "0.0₃0899"
$ws.Range("E30").Value = "  -8.72%  "
$ws.Range("D31").Value = "'505.97"
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("D32").Value = "'7.80"
$ws.Range("E32").Value = "  -6.92%  "
$ws.Range("E33").Value = "  -7.28%  "
$ws.Range("E34").Value = "  -8.36%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'158.93"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -9.55%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'18.48"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'18.57"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  -7.40%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'1.67"
$ws.Range("E42").Value = "  -7.51%  "
$ws.Range("E43").Value = "  -7.36%  "
$ws.Range("D44").Value = "'4.74"
$ws.Range("E44").Value = "  -8.57%  "
$ws.Range("D45").Value = "'38.57"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("D46").Value = "'2.29"
$ws.Range("E46").Value = "  -8.98%  "
$ws.Range("D47").Value = "'141.44"
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("D48").Value = "'3.48"
$ws.Range("E48").Value = "  -6.51%  "
$ws.Range("E49").Value = "  -7.85%  "
$ws.Range("D50").Value = "'1.58"
$ws.Range("E50").Value = "  -8.51%  "
$ws.Range("E51").Value = "  -2.78%  "
